$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose year-over-year trajectory (columns K..AS) must be flattened to
# the value held in column J (the "year 0" value) for that row.
$rows = @(3, 4, 5, 6, 9)

# Column J is index 10; columns K..AS are indices 11..45.
$jCol = 10
$startCol = 11
$endCol = 45

foreach ($r in $rows) {
    $baseValue = $ws.Cells.Item($r, $jCol).Value()
    for ($c = $startCol; $c -le $endCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $baseValue
    }
}
